# Inserção de dados de temperatura de plantas
# Adds a new worksheet "Planilha1" with plant temperature data.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Planilha1"

# Header row (row 2)
$ws.Cells.Item(2, 1).Value = "IdPlanta"
$ws.Cells.Item(2, 2).Value = "Planta"
$ws.Cells.Item(2, 3).Value = "Tipo"
$ws.Cells.Item(2, 4).Value = "Temperatura Miníma (ºC)"
$ws.Cells.Item(2, 5).Value = "Temperatura Maxíma(ºC)"
$ws.Cells.Item(2, 6).Value = "Temp. Min. Ideal (ºC)"
$ws.Cells.Item(2, 7).Value = "Temp. Max. Ideal (ºC)"
$ws.Cells.Item(2, 8).Value = "Tempo médio de germinação (Dias) / (vezes ao ano)"

# Data rows
$data = @(
    @(1, "Tomate",     "Fruta",   10, 34, 15, 25, 90),
    @(2, "Batata",     "legume",  10, 22, 15, 20, 130),
    @(3, "Cenoura",    "Legume",  7,  30, 16, 22, 90),
    @(4, "Beterraba",  "Legume",  10, 24, 10, 20, 85),
    @(5, "Abobrinha",  "Fruta",   15, 35, 18, 27, 60),
    @(6, "Alface",     "Verdura", 7,  24, 15, 19, 55),
    @(7, "Acelga",     "Verdura", 5,  30, 15, 19, 70),
    @(8, "Brocólis",   "Verdura", 6,  28, 21, 25, 90),
    @(9, "Chuchu",     "Verdura", 13, 27, 15, 25, 100),
    @(10, "Couve-flor","Verdura", 14, 25, 15, 20, 115),
    @(11, "Orquídeas", "Flor",    15, 35, 21, 28, 2)
)

$r = 3
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# Select the new sheet and set active cell like original file
$ws.Activate()
$ws.Range("F19").Select()
